# Title slide (slide 1): split the title run
#   "ALTO New Transport using HTTP/2"
# into two runs:
#   "ALTO over " + "New Transport"
# (both keep the same sz/color/font formatting).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$oldTitle = "ALTO New Transport using HTTP/2"
$newTitle = "ALTO over New Transport"

$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldTitle) + 1

# Replace the whole title text first (keeps it as a single run for now).
$titleRange = $tr.Characters($startPos, $oldTitle.Length)
$titleRange.Text = $newTitle

# Re-fetch the text and locate the "New Transport" portion so we can split
# it into its own run by touching its formatting (forces the engine to
# break the run at that boundary while preserving the existing sz/color/
# font attributes already present on the run).
$fullText = $tr.Text
$splitWord = "New Transport"
$splitPos = $fullText.IndexOf($splitWord, $startPos - 1) + 1

$secondRun = $tr.Characters($splitPos, $splitWord.Length)
$secondRun.Font.Size = $secondRun.Font.Size
